# Weekly data refresh: insert this week's new price records (date 44551) at
# the top of the Zanahoria / Vega Monumental Concepción data block, shifting
# the existing history down by two rows (one "Primera" row + one "Segunda"
# row per week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (currently starting at row 132) down by two rows so
# we can insert the new week's two quality-grade records (Primera/Segunda).
$ws.Rows("132:133").Insert()

# New row 132: Zanahoria, "Primera" quality, week of 44551.
$ws.Range("A132").Value = 11
$ws.Range("B132").Value = "Vega Monumental Concepción"
$ws.Range("C132").Value = "Bíobío"
$ws.Range("D132").Value = 44551
$ws.Range("E132").Value = 8
$ws.Range("F132").Value = 100114013
$ws.Range("G132").Value = "Zanahoria"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 6500
$ws.Range("L132").Value = 7000
$ws.Range("M132").Value = 6700
$ws.Range("N132").Value = "$/saco 20 kilos"
$ws.Range("O132").Value = "Región Metropolitana"
$ws.Range("P132").Value = 335
$ws.Range("Q132").Value = 20
$ws.Range("R132").Value = "Hortaliza"

# New row 133: Zanahoria, "Segunda" quality, same week (44551).
$ws.Range("A133").Value = 11
$ws.Range("B133").Value = "Vega Monumental Concepción"
$ws.Range("C133").Value = "Bíobío"
$ws.Range("D133").Value = 44551
$ws.Range("E133").Value = 8
$ws.Range("F133").Value = 100114013
$ws.Range("G133").Value = "Zanahoria"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Segunda"
$ws.Range("J133").Value = 300
$ws.Range("K133").Value = 5500
$ws.Range("L133").Value = 5500
$ws.Range("M133").Value = 5500
$ws.Range("N133").Value = "$/saco 20 kilos"
$ws.Range("O133").Value = "Región Metropolitana"
$ws.Range("P133").Value = 275
$ws.Range("Q133").Value = 20
$ws.Range("R133").Value = "Hortaliza"
